$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Nysteikt vaffel"
$ws.Range("A2").Select()
